$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at 376-377, pushing the existing rows
# (previously 376-492) down to 378-494.
$ws.Rows("376:377").Insert()

# Row 376: new "Especial" quality entry for date 44932 (2023-01-06)
$ws.Range("A376").Value = 4
$ws.Range("B376").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C376").Value = "Los Lagos"
$ws.Range("D376").Value = 44932
$ws.Range("E376").Value = 10
$ws.Range("F376").Value = "Fruta"
$ws.Range("G376").Value = 100102
$ws.Range("H376").Value = "Cítricos"
$ws.Range("I376").Value = 100102006
$ws.Range("J376").Value = "Pomelo"
$ws.Range("K376").Value = "Start Ruby"
$ws.Range("L376").Value = "Especial"
$ws.Range("M376").Value = 80
$ws.Range("N376").Value = 15000
$ws.Range("O376").Value = 15000
$ws.Range("P376").Value = 15000
$ws.Range("Q376").Value = "$/caja 14 kilos empedrada"
$ws.Range("R376").Value = "Región de O'Higgins"
$ws.Range("S376").Value = 1071
$ws.Range("T376").Value = 14

# Row 377: new "Primera" quality entry for the same date 44932
$ws.Range("A377").Value = 4
$ws.Range("B377").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C377").Value = "Los Lagos"
$ws.Range("D377").Value = 44932
$ws.Range("E377").Value = 10
$ws.Range("F377").Value = "Fruta"
$ws.Range("G377").Value = 100102
$ws.Range("H377").Value = "Cítricos"
$ws.Range("I377").Value = 100102006
$ws.Range("J377").Value = "Pomelo"
$ws.Range("K377").Value = "Start Ruby"
$ws.Range("L377").Value = "Primera"
$ws.Range("M377").Value = 80
$ws.Range("N377").Value = 12000
$ws.Range("O377").Value = 12000
$ws.Range("P377").Value = 12000
$ws.Range("Q377").Value = "$/caja 14 kilos empedrada"
$ws.Range("R377").Value = "Región de O'Higgins"
$ws.Range("S377").Value = 857
$ws.Range("T377").Value = 14

# Match the date-column number formatting used by the rest of column D.
$ws.Range("D376:D377").NumberFormat = $ws.Range("D378").NumberFormat
